# Update the "Förändrad" (Changed) date column (C) for rows 2-11
# from 2023-10-05 (serial 45204) to 2023-10-08 (serial 45207).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45207
    }
}
